# Normalize the "Recorded By" (column G) values on the "Session Analysis
# Results" sheet: the list of recorders in each cell is re-ordered into
# ordinal/ASCII sort order (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com").
#
# Only three distinct combinations in this workbook are out of order; all
# other already-sorted combinations (e.g. "System, backup@backdoor.com",
# or single recorders) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colG = 7
$lastRow = $ws.UsedRange.Rows.Count

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
        $changed++
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, backup@backdoor.com, system"
        $changed++
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
        $changed++
    }
}

Write-Host "Reordered 'Recorded By' values in $changed cell(s)."
